$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = [string]"65.415.27"
$ws.Cells.Item(2,5).Value = [string]"  +0.18%  "
$ws.Cells.Item(3,4).Value = [string]"3.539.96"
$ws.Cells.Item(3,5).Value = [string]"  +3.36%  "
$ws.Cells.Item(4,5).Value = [string]"  -0.11%  "
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "597.51"
$ws.Cells.Item(5,5).Value = [string]"  +0.60%  "
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "139.55"
$ws.Cells.Item(6,5).Value = [string]"  +4.11%  "
$ws.Cells.Item(7,4).Value = [string]"3.539.20"
$ws.Cells.Item(7,5).Value = [string]"  +3.39%  "
$ws.Cells.Item(8,5).Value = [string]"  +0.14%  "
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = "0.496"
$ws.Cells.Item(9,5).Value = [string]"  +1.59%  "
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = "0.126"
$ws.Cells.Item(10,5).Value = [string]"  +3.66%  "
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "7.21"
$ws.Cells.Item(11,5).Value = [string]"  -3.37%  "
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "0.392"
$ws.Cells.Item(12,5).Value = [string]"  +4.47%  "
$ws.Cells.Item(13,4).Value = [string]"4.142.43"
$ws.Cells.Item(13,5).Value = [string]"  +3.51%  "
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "0.0000188"
$ws.Cells.Item(14,5).Value = [string]"  +5.25%  "
$ws.Cells.Item(15,2).Value = [string]"Avalanche"
$ws.Cells.Item(15,3).Value = [string]"https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "26.89"
$ws.Cells.Item(15,5).Value = [string]"  +2.47%  "
$ws.Cells.Item(16,2).Value = [string]"WrappedEther"
$ws.Cells.Item(16,3).Value = [string]"https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16,4).Value = [string]"3.531.99"
$ws.Cells.Item(16,5).Value = [string]"  +1.47%  "
$ws.Cells.Item(17,5).Value = [string]"  +1.62%  "
$ws.Cells.Item(18,4).Value = [string]"65.246.32"
$ws.Cells.Item(18,5).Value = [string]"  -0.09%  "
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "10.30"
$ws.Cells.Item(19,5).Value = [string]"  +4.12%  "
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = "5.86"
$ws.Cells.Item(20,5).Value = [string]"  +2.76%  "
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = "14.24"
$ws.Cells.Item(21,5).Value = [string]"  +4.72%  "
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "396.30"
$ws.Cells.Item(22,5).Value = [string]"  +1.31%  "
$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = "0.571"
$ws.Cells.Item(23,5).Value = [string]"  +5.60%  "
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = "74.54"
$ws.Cells.Item(24,5).Value = [string]"  +2.05%  "
$ws.Cells.Item(25,4).Value = [string]"3.685.16"
$ws.Cells.Item(25,5).Value = [string]"  +3.28%  "
$ws.Cells.Item(26,5).Value = [string]"  +0.10%  "
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "0.0000115"
$ws.Cells.Item(27,5).Value = [string]"  +9.15%  "
$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value = "7.75"
$ws.Cells.Item(28,5).Value = [string]"  +9.09%  "
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = "0.997"
$ws.Cells.Item(29,5).Value = [string]"  -0.26%  "
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = "2.27"
$ws.Cells.Item(30,5).Value = [string]"  +1.22%  "
$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = "8.27"
$ws.Cells.Item(31,5).Value = [string]"  +1.73%  "
$ws.Cells.Item(32,4).Value = [string]"3.560.16"
$ws.Cells.Item(32,5).Value = [string]"  +3.76%  "
$ws.Cells.Item(33,5).Value = [string]"  +0.04%  "
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "23.91"
$ws.Cells.Item(34,5).Value = [string]"  +5.96%  "
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "0.146"
$ws.Cells.Item(35,5).Value = [string]"  +0.61%  "
$ws.Cells.Item(36,5).Value = [string]"  +1.27%  "
$ws.Cells.Item(37,2).Value = [string]"Aptos"
$ws.Cells.Item(37,3).Value = [string]"https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value = "7.00"
$ws.Cells.Item(37,5).Value = [string]"  +3.39%  "
$ws.Cells.Item(38,2).Value = [string]"Monero"
$ws.Cells.Item(38,3).Value = [string]"https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value = "168.96"
$ws.Cells.Item(38,5).Value = [string]"  -1.49%  "
$ws.Cells.Item(39,5).Value = [string]"  +2.36%  "
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "4.93"
$ws.Cells.Item(40,5).Value = [string]"  +2.75%  "
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = "0.0804"
$ws.Cells.Item(41,5).Value = [string]"  +5.11%  "
$ws.Cells.Item(42,5).Value = [string]"  +1.52%  "
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "26.69"
$ws.Cells.Item(43,5).Value = [string]"  +22.10%  "
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = "42.67"
$ws.Cells.Item(44,5).Value = [string]"  -1.62%  "
$ws.Cells.Item(45,5).Value = [string]"  -0.14%  "
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = "4.45"
$ws.Cells.Item(46,5).Value = [string]"  +1.74%  "
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = "1.19"
$ws.Cells.Item(47,5).Value = [string]"  +10.02%  "
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "1.67"
$ws.Cells.Item(48,5).Value = [string]"  +4.20%  "
$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value = "6.83"
$ws.Cells.Item(49,5).Value = [string]"  +4.98%  "
$ws.Cells.Item(50,4).Value = [string]"2.368.62"
$ws.Cells.Item(50,5).Value = [string]"  +8.27%  "
$ws.Cells.Item(51,2).Value = [string]"LidoDAOToken"
$ws.Cells.Item(51,3).Value = [string]"https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "2.37"
$ws.Cells.Item(51,5).Value = [string]"  +26.24%  "
